# Generate Report for Handback
#
# The handback for file "9a7ef578-c2fd-423a-8afd-4c24e7916674.md" failed:
# the handback file name did not match the handoff file name. Update the
# status for that row on the Overview sheet and on each locale sheet
# (zh-cn / de-de), and record the error detail in the "Error Detail"
# column for each locale sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$errorDetailZhCn = "Handback file name: d2agn3fi.2ht is different with handoff file name: 9a7ef578-c2fd-423a-8afd-4c24e7916674.8061c377b658f4820c7373f3916348e8b4bdd5f9.zh-cn."
$errorDetailDeDe = "Handback file name: d2agn3fi.2ht is different with handoff file name: 9a7ef578-c2fd-423a-8afd-4c24e7916674.8061c377b658f4820c7373f3916348e8b4bdd5f9.de-de."

# --- Overview sheet: update the status for the second file in both locale columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet: update status + error detail ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K3").Value = $errorDetailZhCn

# --- de-de sheet: update status + error detail ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("K3").Value = $errorDetailDeDe
